$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 116
$ws.Range("H116").Value = 8811.299999999999
$ws.Range("I116").Value = 6931.5557
$ws.Range("J116").Value = 10349.272
$ws.Range("K116").Value = 6931.5557
$ws.Range("L116").Value = 10349.272
$ws.Range("M116").Value = -3489.5557
$ws.Range("N116").Value = -17233.272

# ALC row 132
$ws.Range("H132").Value = 2502.7144
$ws.Range("I132").Value = 2702.8
$ws.Range("J132").Value = 2002.5
$ws.Range("K132").Value = 8108.400000000001
$ws.Range("L132").Value = 6007.5
$ws.Range("M132").Value = -5578.400000000001

# ALC row 134
$ws.Range("H134").Value = 45780
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 45780
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 45780
$ws.Range("N134").Value = -55920

# ALC row 138
$ws.Range("H138").Value = 4156.778
$ws.Range("I138").Value = 3313
$ws.Range("J138").Value = 4940.2856
$ws.Range("K138").Value = 9939
$ws.Range("L138").Value = 14820.8568
$ws.Range("M138").Value = -4799

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45
$ws.Range("H45").Value = 8890.777
$ws.Range("I45").Value = 11365.083
$ws.Range("J45").Value = 3942.1667
$ws.Range("K45").Value = 11365.083
$ws.Range("L45").Value = 3942.1667
$ws.Range("M45").Value = -10988.083
$ws.Range("N45").Value = -4696.1667

# ARM row 61
$ws.Range("H61").Value = 1012
$ws.Range("I61").Value = 1012
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1012
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -800
$ws.Range("N61").ClearContents()

# ARM row 74
$ws.Range("H74").Value = 5328.1665
$ws.Range("I74").Value = 2950
$ws.Range("J74").Value = 6517.25
$ws.Range("K74").Value = 2950
$ws.Range("L74").Value = 6517.25
$ws.Range("M74").Value = -2076
$ws.Range("N74").Value = -8265.25

# ARM row 77
$ws.Range("H77").Value = 5328.1665
$ws.Range("I77").Value = 2950
$ws.Range("J77").Value = 6517.25
$ws.Range("K77").Value = 14750
$ws.Range("L77").Value = 32586.25
$ws.Range("M77").Value = -10382
$ws.Range("N77").Value = -41322.25

# ARM row 136
$ws.Range("H136").Value = 1012
$ws.Range("I136").Value = 1012
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3036
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -486
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws.Range("H105").Value = 2165.1956
$ws.Range("I105").Value = 1946.5862
$ws.Range("J105").Value = 2538.1177
$ws.Range("K105").Value = 1946.5862
$ws.Range("L105").Value = 2538.1177
$ws.Range("M105").Value = -199.5862
$ws.Range("N105").Value = -6032.1177

# BSM row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 2475.577
$ws.Range("I31").Value = 1548
$ws.Range("J31").Value = 3403.1538
$ws.Range("K31").Value = 1548
$ws.Range("L31").Value = 3403.1538
$ws.Range("M31").Value = -1253

# CRP row 34
$ws.Range("H34").Value = 2475.577
$ws.Range("I34").Value = 1548
$ws.Range("J34").Value = 3403.1538
$ws.Range("K34").Value = 1548
$ws.Range("L34").Value = 3403.1538
$ws.Range("M34").Value = -1346

# CRP row 50
$ws.Range("H50").Value = 29916.545
$ws.Range("I50").Value = 25082
$ws.Range("J50").Value = 30400
$ws.Range("K50").Value = 25082
$ws.Range("L50").Value = 30400
$ws.Range("M50").Value = -24457
$ws.Range("N50").Value = -31650

# CRP row 86
$ws.Range("H86").Value = 4945.0557
$ws.Range("I86").Value = 4937.2856
$ws.Range("J86").Value = 4972.25
$ws.Range("K86").Value = 4937.2856
$ws.Range("L86").Value = 4972.25
$ws.Range("M86").Value = -3814.2856
$ws.Range("N86").Value = -7218.25

# CRP row 89
$ws.Range("H89").Value = 4945.0557
$ws.Range("I89").Value = 4937.2856
$ws.Range("J89").Value = 4972.25
$ws.Range("K89").Value = 24686.428
$ws.Range("L89").Value = 24861.25
$ws.Range("M89").Value = -19070.428
$ws.Range("N89").Value = -36093.25

# CRP row 134
$ws.Range("H134").Value = 4269.9414
$ws.Range("I134").Value = 4776.231
$ws.Range("J134").Value = 2624.5
$ws.Range("K134").Value = 14328.693
$ws.Range("L134").Value = 7873.5
$ws.Range("M134").Value = -11793.693
$ws.Range("N134").Value = -12943.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 8
$ws.Range("H8").Value = 654.3333
$ws.Range("I8").Value = 654.3333
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1962.9999
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1823.9999

# CUL row 86
$ws.Range("H86").Value = 1930
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1930
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 5790
$ws.Range("N86").Value = -8162
$ws.Range("M86").ClearContents()

# CUL row 89
$ws.Range("H89").Value = 1930
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1930
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 17370
$ws.Range("N89").Value = -29226
$ws.Range("M89").ClearContents()

# CUL row 103
$ws.Range("H103").Value = 681
$ws.Range("I103").Value = 681
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 2043
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -1164
$ws.Range("N103").ClearContents()

# CUL row 107
$ws.Range("H107").Value = 1071.5
$ws.Range("I107").Value = 627.6
$ws.Range("J107").Value = 1318.1111
$ws.Range("K107").Value = 1882.8
$ws.Range("L107").Value = 3954.3333
$ws.Range("M107").Value = 37.19999999999982
$ws.Range("N107").Value = -7794.3333

# CUL row 119
$ws.Range("H119").Value = 2840
$ws.Range("I119").Value = 2840
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 8520
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -3682

# CUL row 131
$ws.Range("H131").Value = 1373.0182
$ws.Range("I131").Value = 413.72726
$ws.Range("J131").Value = 1612.841
$ws.Range("K131").Value = 1241.18178
$ws.Range("L131").Value = 4838.522999999999
$ws.Range("M131").Value = 3798.81822
$ws.Range("N131").Value = -14918.523

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 13738.25
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 13738.25
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 13738.25
$ws.Range("N80").Value = -15734.25

# GSM row 83
$ws.Range("H83").Value = 13738.25
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 13738.25
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 68691.25
$ws.Range("N83").Value = -78675.25

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 2351
$ws.Range("I16").Value = 3256.7
$ws.Range("J16").Value = 1057.1428
$ws.Range("K16").Value = 3256.7
$ws.Range("L16").Value = 1057.1428
$ws.Range("M16").Value = -3086.7
$ws.Range("N16").Value = -1397.1428

# LTW row 43
$ws.Range("H43").Value = 300181.56
$ws.Range("I43").Value = 382863.16
$ws.Range("J43").Value = 217500
$ws.Range("K43").Value = 382863.16
$ws.Range("L43").Value = 217500
$ws.Range("M43").Value = -382670.16
$ws.Range("N43").Value = -217886

# LTW row 55
$ws.Range("H55").Value = 958.2
$ws.Range("I55").Value = 1155.25
$ws.Range("J55").Value = 170
$ws.Range("K55").Value = 1155.25
$ws.Range("L55").Value = 170
$ws.Range("M55").Value = -982.25
$ws.Range("N55").Value = -516

# LTW row 82
$ws.Range("H82").Value = 800.125
$ws.Range("I82").Value = 812.25
$ws.Range("J82").Value = 788
$ws.Range("K82").Value = 812.25
$ws.Range("L82").Value = 788
$ws.Range("M82").Value = -451.25
$ws.Range("N82").Value = -1510

# LTW row 85
$ws.Range("H85").Value = 800.125
$ws.Range("I85").Value = 812.25
$ws.Range("J85").Value = 788
$ws.Range("K85").Value = 812.25
$ws.Range("L85").Value = 788
$ws.Range("M85").Value = 435.75
$ws.Range("N85").Value = -3284

# LTW row 93
$ws.Range("H93").Value = 2015.5
$ws.Range("I93").Value = 2112.3157
$ws.Range("J93").Value = 1647.6
$ws.Range("K93").Value = 2112.3157
$ws.Range("L93").Value = 1647.6
$ws.Range("M93").Value = -864.3157000000001

# LTW row 122
$ws.Range("H122").Value = 3006.75
$ws.Range("I122").Value = 2640.2307
$ws.Range("J122").Value = 3687.4285
$ws.Range("K122").Value = 7920.6921
$ws.Range("L122").Value = 11062.2855
$ws.Range("M122").Value = -5470.6921
$ws.Range("N122").Value = -15962.2855

# LTW row 123
$ws.Range("H123").Value = 43429
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 43429
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 43429
$ws.Range("N123").Value = -53229

$ws = $wb.Worksheets.Item("WVR")
# WVR row 112
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# WVR row 113
$ws.Range("H113").Value = 599.6667
$ws.Range("I113").Value = 599.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1799.0001
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 370.9999

# WVR row 122
$ws.Range("H122").Value = 3475.3157
$ws.Range("I122").Value = 3359.7778
$ws.Range("J122").Value = 5555
$ws.Range("K122").Value = 10079.3334
$ws.Range("L122").Value = 16665
$ws.Range("M122").Value = -7629.3334

# WVR row 126
$ws.Range("H126").Value = 4934.8887
$ws.Range("I126").Value = 4997.5415
$ws.Range("J126").Value = 4809.5835
$ws.Range("K126").Value = 14992.6245
$ws.Range("L126").Value = 14428.7505
$ws.Range("M126").Value = -12522.6245
$ws.Range("N126").Value = -19368.7505

# WVR row 139
$ws.Range("H139").Value = 47048.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 47048.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 47048.332
$ws.Range("N139").Value = -57328.332
